$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9973810911178589
$ws.Range("B1").Value = 2.10918927192688
$ws.Range("C1").Value = 7.018016815185547
$ws.Range("D1").Value = 2.098139524459839
$ws.Range("E1").Value = 1.376373052597046
